$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.197.91"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "4.049.66"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("D5").Value = "'573.92"
$ws.Range("E5").Value = "  +7.78%  "
$ws.Range("D6").Value = "'152.10"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "4.046.42"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("D8").Value = "'0.695"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("E10").Value = "  +2.45%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "'53.94"
$ws.Range("E12").Value = "  +13.43%  "
$ws.Range("D13").Value = "'0.0000327"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "'11.20"
$ws.Range("E14").Value = "  +4.70%  "
$ws.Range("D15").Value = "4.696.48"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "4.057.53"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "'14.37"
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "'20.91"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("E19").Value = "  +2.99%  "
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "73.212.76"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").Value = "'444.97"
$ws.Range("E22").Value = "  +4.41%  "
$ws.Range("D23").Value = "'98.52"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").Value = "'4.52"
$ws.Range("E24").Value = "  +7.17%  "
$ws.Range("D25").Value = "'3.57"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("D26").Value = "'14.68"
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("E27").Value = "  +18.98%  "
$ws.Range("D28").Value = "'11.42"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").Value = "'5.96"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").Value = "'37.17"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("E32").Value = "  +10.48%  "
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("D34").Value = "'13.73"
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("D35").Value = "'695.56"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").Value = "'48.60"
$ws.Range("E36").Value = "  +11.07%  "
$ws.Range("D37").Value = "'68.77"
$ws.Range("E37").Value = "  +5.20%  "
$ws.Range("D38").Value = "0.0₃0905"
$ws.Range("E38").Value = "  +9.70%  "
$ws.Range("E39").Value = "  +4.12%  "
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  +1.61%  "
$ws.Range("D45").Value = "'0.0496"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("D48").Value = "'2.75"
$ws.Range("E48").Value = "  +5.05%  "
$ws.Range("D49").Value = "'2.23"
$ws.Range("E49").Value = "  +11.61%  "
$ws.Range("D50").Value = "'3.51"
$ws.Range("E50").Value = "  +6.56%  "
$ws.Range("D51").Value = "'3.06"
$ws.Range("E51").Value = "  +2.62%  "
